$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-18 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-19 Friday", 2)

$d.Content.Find.Execute("81×37=", $true, $false, $false, $false, $false, $true, 1, $false, "42×82=", 2)
$d.Content.Find.Execute("85×45=", $true, $false, $false, $false, $false, $true, 1, $false, "35×73=", 2)
$d.Content.Find.Execute("26×61=", $true, $false, $false, $false, $false, $true, 1, $false, "90×52=", 2)
$d.Content.Find.Execute("35×91=", $true, $false, $false, $false, $false, $true, 1, $false, "36×74=", 2)
$d.Content.Find.Execute("39×88=", $true, $false, $false, $false, $false, $true, 1, $false, "23×62=", 2)

$d.Content.Find.Execute("51×55=", $true, $false, $false, $false, $false, $true, 1, $false, "81×25=", 2)
$d.Content.Find.Execute("76×94=", $true, $false, $false, $false, $false, $true, 1, $false, "95×67=", 2)
$d.Content.Find.Execute("39×16=", $true, $false, $false, $false, $false, $true, 1, $false, "44×84=", 2)
$d.Content.Find.Execute("35×36=", $true, $false, $false, $false, $false, $true, 1, $false, "74×79=", 2)
$d.Content.Find.Execute("91×34=", $true, $false, $false, $false, $false, $true, 1, $false, "56×55=", 2)

$d.Content.Find.Execute("48×26=", $true, $false, $false, $false, $false, $true, 1, $false, "59×91=", 2)
$d.Content.Find.Execute("18×70=", $true, $false, $false, $false, $false, $true, 1, $false, "87×91=", 2)
$d.Content.Find.Execute("52×78=", $true, $false, $false, $false, $false, $true, 1, $false, "94×18=", 2)
$d.Content.Find.Execute("92×30=", $true, $false, $false, $false, $false, $true, 1, $false, "26×23=", 2)
$d.Content.Find.Execute("51×39=", $true, $false, $false, $false, $false, $true, 1, $false, "13×12=", 2)

$d.Content.Find.Execute("43×67=", $true, $false, $false, $false, $false, $true, 1, $false, "96×26=", 2)
$d.Content.Find.Execute("95×74=", $true, $false, $false, $false, $false, $true, 1, $false, "30×93=", 2)
$d.Content.Find.Execute("64×89=", $true, $false, $false, $false, $false, $true, 1, $false, "72×97=", 2)
$d.Content.Find.Execute("26×90=", $true, $false, $false, $false, $false, $true, 1, $false, "12×49=", 2)
$d.Content.Find.Execute("96×99=", $true, $false, $false, $false, $false, $true, 1, $false, "55×37=", 2)

$d.Content.Find.Execute("96×23=", $true, $false, $false, $false, $false, $true, 1, $false, "38×93=", 2)
$d.Content.Find.Execute("21×79=", $true, $false, $false, $false, $false, $true, 1, $false, "14×94=", 2)
$d.Content.Find.Execute("72×73=", $true, $false, $false, $false, $false, $true, 1, $false, "29×73=", 2)
$d.Content.Find.Execute("69×94=", $true, $false, $false, $false, $false, $true, 1, $false, "58×70=", 2)
$d.Content.Find.Execute("89×52=", $true, $false, $false, $false, $false, $true, 1, $false, "89×37=", 2)
